# Rename the sheets (new task order IDs) and update the stim/file-name
# values in column B of each sheet to reflect the new practice/final
# task-order run.

$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16509961086022134"
$wb.Worksheets.Item(2).Name = "NB_TO-16509961113303747"
$wb.Worksheets.Item(3).Name = "RS_TO-16509961113303747"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509961113783798"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650996111442379"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650996108570213.csv"
$ws1.Range("B3").Value = "GNG_stims-1650996108586251.csv"
$ws1.Range("B4").Value = "go_stims-1650996108586251.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961086022134.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_9-16509961087622132.csv"
$ws2.Range("B3").Value = "OB-16509961094023752.csv"
$ws2.Range("B4").Value = "TB-16509961113144088.csv"
$ws2.Range("B5").Value = "ZB-match_7-16509961090582466.csv"
$ws2.Range("B6").Value = "TB-16509961104984162.csv"
$ws2.Range("B7").Value = "ZB-match_1-16509961087382088.csv"
$ws2.Range("B8").Value = "OB-16509961095383804.csv"
$ws2.Range("B9").Value = "TB-16509961105383778.csv"
$ws2.Range("B10").Value = "OB-16509961096423776.csv"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509961113463814.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961113303747.csv"
$ws4.Range("B4").Value = "MM_stims-16509961113624089.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961113463814.csv"
$ws4.Range("B6").Value = "MM_stims-16509961113783798.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961113624089.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16509961113783798.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961113944252.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961114104087.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961114264095.csv"
